# "Control Stock Web" sheet: drop the "10962736022" item (old row 2) and
# re-sort the remaining three stock rows into their new order:
#   old row 5 -> new row 2
#   old row 4 -> new row 3
#   old row 3 -> new row 4
#
# Cell values such as "10962389016" and "$ 165.990" must stay as literal
# text (matching how the file was originally authored), so we move the
# existing rows with Range.Copy instead of re-typing the strings through
# Range.Value - the latter makes Excel "smart" re-parse numeric-looking
# text into real numbers/currency, which is not what we want here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage the four current data rows (2-5) in a scratch area, far away from
# the used range, so copying them into their destinations can't clobber a
# row that hasn't been read yet.
$ws.Range("A2:D2").Copy($ws.Range("A20:D20"))  # 10962736022 (dropped)
$ws.Range("A3:D3").Copy($ws.Range("A21:D21"))  # 10930745010
$ws.Range("A4:D4").Copy($ws.Range("A22:D22"))  # 10962389018
$ws.Range("A5:D5").Copy($ws.Range("A23:D23"))  # 10962389016

# Write the new rows 2-4 from the staged copies, in the target order.
$ws.Range("A23:D23").Copy($ws.Range("A2:D2"))  # 10962389016 -> row 2
$ws.Range("A22:D22").Copy($ws.Range("A3:D3"))  # 10962389018 -> row 3
$ws.Range("A21:D21").Copy($ws.Range("A4:D4"))  # 10930745010 -> row 4

# Clean up the scratch area, then drop the now-redundant old row 5
# (dimension shrinks from A1:D5 to A1:D4).
$ws.Range("A20:D23").Clear()
$ws.Rows.Item(5).Delete()
